# 28-08-2025(scroll view fixed )
#
# The "Total Amount" column (S) is no longer needed, so remove it outright
# (entire-column delete), which shifts every column from T onward one slot
# to the left and drops the now-unused "Total Amount" shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S1").EntireColumn.Delete()

# Restore a sensible view: scroll the grid over to the right-hand columns
# and leave the selection where the user left it after the edit.
[void]$ws.Range("W27").Select()
